$wb = $excel.ActiveWorkbook

function Set-Row($SheetName, $Row, $Values) {
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($col in $Values.Keys) {
        $addr = "$col$Row"
        $ws.Range($addr).Value = $Values[$col]
    }
}

function Clear-Cell($SheetName, $Addr) {
    $ws = $wb.Worksheets.Item($SheetName)
    $ws.Range($Addr).ClearContents()
}

# ===== ALC =====
Set-Row "ALC" 34 @{
    H = 11514.667; I = 1817.6; J = 60000; K = 1817.6; L = 60000; M = -1614.6; N = -60406
}
Set-Row "ALC" 36 @{
    H = 11514.667; I = 1817.6; J = 60000; K = 1817.6; L = 60000; M = -1102.6; N = -61430
}
Set-Row "ALC" 40 @{
    H = 1602.579; I = 1511; J = 1685; K = 1511; L = 1685; M = -1336; N = -2035
}
Set-Row "ALC" 129 @{
    H = 1210.78; I = 399.4; J = 1300.9333; K = 1198.2; L = 3902.7999; M = 3801.8; N = -13902.7999
}
Set-Row "ALC" 138 @{
    H = 236024.66; I = 1468.2285; J = 609182.5600000001; K = 4404.6855; L = 1827547.68; M = 735.3145000000004; N = -1837827.68
}

# ===== ARM =====
Set-Row "ARM" 123 @{
    H = 543596; J = 543596; L = 543596; N = -553396
}

# ===== CRP =====
Set-Row "CRP" 23 @{
    H = 0; I = 0; K = 0
}
Clear-Cell "CRP" "M23"

Set-Row "CRP" 27 @{
    H = 0; I = 0; K = 0
}
Clear-Cell "CRP" "M27"

Set-Row "CRP" 31 @{
    H = 1998.7567; I = 1291.2903; J = 2508.7908; K = 1291.2903; L = 2508.7908; M = -996.2902999999999; N = -3098.7908
}
Set-Row "CRP" 34 @{
    H = 1998.7567; I = 1291.2903; J = 2508.7908; K = 1291.2903; L = 2508.7908; M = -1089.2903; N = -2912.7908
}
Set-Row "CRP" 99 @{
    H = 1424.0714; I = 1311.1818; J = 1838; K = 1311.1818; L = 1838; M = 186.8181999999999; N = -4834
}
Set-Row "CRP" 107 @{
    H = 1231.8518; I = 1287.5294; J = 1137.2; K = 1287.5294; L = 1137.2; M = 632.4706000000001; N = -4977.2
}
Set-Row "CRP" 126 @{
    H = 1424.0714; I = 1311.1818; J = 1838; K = 3933.5454; L = 5514; M = -1463.5454; N = -10454
}

# ===== CUL =====
Set-Row "CUL" 68 @{
    H = 1511.0117; I = 1226.3334; J = 1831.275; K = 3679.0002; L = 5493.825000000001; M = -2868.0002; N = -7115.825000000001
}
Set-Row "CUL" 71 @{
    H = 1511.0117; I = 1226.3334; J = 1831.275; K = 11037.0006; L = 16481.475; M = -6981.000599999999; N = -24593.475
}
Set-Row "CUL" 107 @{
    H = 964.5192; I = 697.1724; J = 1301.6086; K = 2091.5172; L = 3904.8258; M = -171.5172000000002; N = -7744.825800000001
}
Set-Row "CUL" 131 @{
    H = 961.8955; I = 471; J = 1028.4576; K = 1413; L = 3085.3728; M = 3627; N = -13165.3728
}
Set-Row "CUL" 132 @{
    H = 631720.25; J = 6556.4287; L = 59007.85830000001; N = -64067.85830000001
}
Set-Row "CUL" 134 @{
    H = 796.25; I = 796.25; K = 2388.75; M = 2681.25
}
Set-Row "CUL" 139 @{
    H = 28725.73; I = 1477.1562; J = 203116.6; K = 4431.4686; L = 609349.8; M = 708.5313999999998; N = -619629.8
}

# ===== GSM =====
Set-Row "GSM" 5 @{
    H = 9995; J = 9995; L = 9995; N = -10219
}
Set-Row "GSM" 123 @{
    H = 16333.333; J = 16333.333; L = 16333.333; N = -21233.333
}

# ===== LTW =====
Set-Row "LTW" 4 @{
    H = 0; I = 0; K = 0
}
Clear-Cell "LTW" "M4"

Set-Row "LTW" 25 @{
    H = 17008; I = 0; J = 17008; K = 0; L = 17008; N = -17468
}
Clear-Cell "LTW" "M25"

Set-Row "LTW" 28 @{
    H = 0; I = 0; K = 0
}
Clear-Cell "LTW" "M28"

Set-Row "LTW" 37 @{
    H = 0; I = 0; K = 0
}
Clear-Cell "LTW" "M37"

Set-Row "LTW" 100 @{
    H = 2007.3572; I = 1418.4546; J = 4166.6665; K = 1418.4546; L = 4166.6665; M = -877.4546; N = -5248.6665
}

# ===== WVR =====
Set-Row "WVR" 24 @{
    H = 70010; J = 70010; L = 70010; N = -70470
}
